$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    1  = 0.22195987269149242
    2  = -0.0059999999463791198
    3  = -0.0039999999524873431
    4  = -0.0079999999135100808
    5  = -0.0029999999516405751
    6  = -0.0019999999472890551
    7  = 0.034240474863218129
    8  = -0.0099999998772188903
    9  = -0.0019999999441777661
    10 = -0.0019999999430986293
    11 = -0.0029999999345005079
    12 = 0.034565185851041846
    13 = -0.0034999999211429156
    14 = -0.0079999998797735117
    15 = -0.00099999993744059879
    16 = -0.0019999999263888846
    17 = -0.001999999921831197
    18 = -0.016554215914565518
    19 = -0.0039999999617847948
    20 = -0.0039999999442983381
    21 = -0.0039999999404622955
    22 = -0.057824049417573242
    23 = -0.0049999999409342522
    24 = -0.019999999806880275
    25 = -0.019999999804234392
    26 = -0.0024999999402677275
    27 = -0.0024999999364059278
    28 = -0.0019999999229316501
    29 = -0.006999999868106066
    30 = -0.059999999408933746
    31 = -0.0069999998581540268
    32 = -0.0099999998318249794
    33 = -0.0039999998829394201
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 1).Value = $values[$row]
}
